$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (removes the stray "mapping-siglas.xlsx" D6 cell)
$ws.Rows("6:6").Delete()

# Re-style D1:D5 and A3:A5 from style 2 to style 1 (copy format from a style-1 cell, e.g. B1)
$ws.Range("B1").Copy()
$ws.Range("D1:D5").PasteSpecial(-4122)
$ws.Range("A3:A5").PasteSpecial(-4122)

# Update cell values (text) for rows 3-5, columns A and D
$ws.Range("A3").Value = "iaest-measure:orden"
$ws.Range("D3").Value = "iaest-measure:siglas"
$ws.Range("A4").Value = "medida"
$ws.Range("D4").Value = "medida"
$ws.Range("A5").Value = "xsd:int"
$ws.Range("D5").Value = "xsd:string"
